$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B23").Value = "Export"
$ws.Range("B28").Value = "überschreiben! Willst du wirklich den Vorgang fortsetzten?"

$ws.Range("G2").Formula = '=SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(D2," ","\u0020"),"!"," \u0021")," ","")'
$ws.Range("G3:G66").Formula = '=SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(D3," ","\u0020"),"!"," \u0021")," ","")'
$ws.Range("G67:G69").Formula = '=SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(D67," ","\u0020"),"!"," \u0021")," ","")'
